$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Ana García Rodríguez
$ws.Range("B2").Value = "ana.garcia@hackless.com"
$ws.Range("D2").Value = "Operadora de Línea de Producción"
$ws.Range("F2").Value = "+54 9 11 1234-5678"

# Row 3 - Carlos López Martínez
$ws.Range("B3").Value = "carlos.lopez@hackless.com"
$ws.Range("D3").Value = "Supervisor de Seguridad Industrial"
$ws.Range("F3").Value = "+54 9 11 2345-6789"

# Row 4 - María José Fernández (name changed)
$ws.Range("A4").Value = "María José Fernández"
$ws.Range("B4").Value = "maria.fernandez@hackless.com"
$ws.Range("D4").Value = "Jefe de Recursos Humanos"
$ws.Range("E4").Value = "Administración"
$ws.Range("F4").Value = "+54 9 11 3456-7890"

# Row 5 - Juan Carlos Sánchez (name changed)
$ws.Range("A5").Value = "Juan Carlos Sánchez"
$ws.Range("B5").Value = "juan.sanchez@hackless.com"
$ws.Range("D5").Value = "Técnico en Mantenimiento Mecánico"
$ws.Range("F5").Value = "+54 9 11 4567-8901"

# Row 6 - Laura Patricia Morales (name changed)
$ws.Range("A6").Value = "Laura Patricia Morales"
$ws.Range("B6").Value = "laura.morales@hackless.com"
$ws.Range("D6").Value = "Auditor Interno de Calidad"
$ws.Range("F6").Value = "+54 9 11 5678-9012"
